$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: column D ("Price") holds numeric-looking text (e.g. "237.54").
# A leading apostrophe forces Excel to keep it stored as text instead of
# silently converting it to a real number (which would also mangle
# trailing zeros like "0.0001500" / "2.200").

# Rows 2-12: simple price refresh, no name/link/volume change
$ws.Range("D2").Value = "'237.54"
$ws.Range("D3").Value = "'21.64"
$ws.Range("D4").Value = "'5.357"
$ws.Range("D5").Value = "'0.05556"
$ws.Range("D7").Value = "'6.453"
$ws.Range("D8").Value = "'0.8025"
$ws.Range("D9").Value = "'1.036"
$ws.Range("D10").Value = "'0.1399"
$ws.Range("D11").Value = "'0.07311"
$ws.Range("D12").Value = "'0.03271"

# Rows 13-27: the coin list shifted up by one (row 13's old coin,
# ProBitToken, wraps around to the bottom at row 27), and each row's
# price/volume-rank text was refreshed to the latest snapshot.
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02872"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09242"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001665"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.257"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04764"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005709"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006258"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.005057"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.001053"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "UpBots"
$ws.Range("C23").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D23").Value = "'0.0004180"
$ws.Range("E23").Value = "22UpBotsUBXT"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'3.969"
$ws.Range("E24").Value = "23LEOLEOBestin24h"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.200"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3289"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("B27").Value = "ProBitToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D27").Value = "'0.1287"
$ws.Range("E27").Value = "26ProBitTokenPROB"

# Rows 40-48: further simple price refreshes (plus one volume-rank
# label that dropped its "Worstin24h" suffix on row 48).
$ws.Range("D40").Value = "'0.04136"
$ws.Range("D41").Value = "'0.007036"
$ws.Range("D42").Value = "'0.003499"
$ws.Range("D43").Value = "'0.1035"
$ws.Range("D44").Value = "'0.008797"
$ws.Range("D47").Value = "'0.6798"
$ws.Range("D48").Value = "'0.03152"
$ws.Range("E48").Value = "47BOLOBOLO"
